$d = $word.ActiveDocument
$d.Content.Find.Execute("denem", $true, $false, $false, $false, $false, $true, 1, $false, "deben", 2)
